$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Mittelstand-Digital Zentrum" value for row 2 (AE2), matching the
# left-aligned style already used by the other AE column entries (e.g. AE3).
$ws.Range("AE2").Value = "Mittelstand-Digital Zentrum Kaiserslautern"
$ws.Range("AE2").HorizontalAlignment = -4131

# Expand the abbreviated "MD Zentrum" label to the full "Mittelstand-Digital Zentrum" name.
$ws.Range("AE3").Value = "Mittelstand-Digital Zentrum Lingen.Münster.Osnabrück"

# Add new test registrations that were used to reproduce/verify the bug fix.
$ws.Range("S11").Value = "Testperson"
$ws.Range("T11").Value = "Testnachname"
$ws.Range("U11").Value = "Testmail@mail.com"
$ws.Range("V11").Value = "TestOrganisation"
$ws.Range("AD11").Value = "Testmail@mail.com"

$ws.Range("S12").Value = "Neuertest"
$ws.Range("T12").Value = "Nachname"
$ws.Range("U12").Value = "Neuemail@mail.de"
$ws.Range("V12").Value = "NeueOrganisation"
$ws.Range("AD12").Value = "Neuemail@mail.de"

# Update the active selection to match the final state of the sheet.
$ws.Range("AA11").Select()
